# Generate Report for Handback
# Update status + error detail for the "855881e7-..." row (row 3) on both
# the zh-cn and de-de localization-status sheets, and widen the
# "Error Detail" column (P) to fit the new message.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Row 3 / "855881e7-3456-493c-be1c-fc7451fe05ac.md" status changes from
# "Ready for handoff" to "Handback transform failed" everywhere it is
# shown: the per-language sheets and the roll-up Overview sheet.
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"

# Populate the (previously empty) Error Detail column for that row.
$zhcn.Range("P3").Value = "Handback file name: n1itf0ta.0lv is different with handoff file name: 855881e7-3456-493c-be1c-fc7451fe05ac.da6129ed12ddd3e6ae2b10a98faad796cf781198.zh-cn."
$dede.Range("P3").Value = "Handback file name: n1itf0ta.0lv is different with handoff file name: 855881e7-3456-493c-be1c-fc7451fe05ac.da6129ed12ddd3e6ae2b10a98faad796cf781198.de-de."

# Widen the Error Detail column so the new message is readable. Excel's
# ColumnWidth (characters) is stored in the XML as ColumnWidth + 5/6, so
# back the input off by that padding to land on an OOXML width of 40.
$targetWidth = 40 - (5 / 6)
$zhcn.Columns.Item(16).ColumnWidth = $targetWidth
$dede.Columns.Item(16).ColumnWidth = $targetWidth
